# Edit script for ManaloK9-Tracker.xlsx
# - Populates rows 2-5 of Sheet1 with booking data (as literal text, matching the
#   original inlineStr/text cell type of the tracker).
# - Adds an AutoFilter over the header row A1:V1 (and the hidden _FilterDatabase
#   defined name Excel creates as a side effect of AutoFilter).
# - Moves the active-cell selection to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 22
$startRow = 2

$values = @(
    "Cancel",
    "Call",
    "2025-01-09",
    "Rodzell Jan Gamboa Cerda",
    "2025-01-15",
    "09",
    "20:28",
    "20:28",
    "2",
    "1",
    "2",
    "1",
    "1",
    "1",
    "2",
    "1,10",
    "1",
    "4",
    "Rodzell Jan Cerda",
    "3000",
    "877",
    "6000",
    "Done",
    "FB",
    "2025-01-01",
    "Rodzell Jan Gamboa Cerda",
    "2025-01-02",
    "09",
    "22:44",
    "22:44",
    "2",
    "1",
    "2",
    "1",
    "1",
    "1",
    "2",
    "1,10",
    "1",
    "4",
    "Rodzell Jan Cerda",
    "3000",
    "877",
    "6000",
    "Reserved",
    "FB",
    "2025-01-01",
    "Rodzell Jan Gamboa Cerda",
    "2025-01-19",
    "09",
    "22:00",
    "10:00",
    "1",
    "1",
    "2",
    "1",
    "1",
    "1",
    "2",
    "1,10",
    "1",
    "4",
    "Rodzell Jan Cerda",
    "3000",
    "877",
    "6000",
    "Done",
    "FB",
    "2025-01-05",
    "Rodzell Jan Gamboa Cerda",
    "2025-01-07",
    "09",
    "23:12",
    "23:08",
    "2",
    "1",
    "2",
    "1",
    "1",
    "1",
    "2",
    "1,10",
    "1",
    "4",
    "Rodzell Jan Cerda",
    "3000",
    "877",
    "6000"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + [Math]::Floor($i / $numCols)
    $col = ($i % $numCols) + 1
    # Leading apostrophe forces text storage (matches the source file, where
    # every data cell - including numbers, dates and times - is stored as text).
    $ws.Cells.Item($row, $col).Value = "'" + $values[$i]
}

# Turn on AutoFilter for the header row; Excel also writes a hidden
# worksheet-scoped _FilterDatabase defined name as a side effect.
[void]$ws.Range("A1:V1").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Sheet1'!`$A`$1:`$V`$1")
$filterName.Visible = $false

# Restore the saved cursor position.
[void]$ws.Range("I15").Select()
